$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 150, pushing existing rows 150-163 down to 152-165.
$ws.Rows.Item(150).Resize(2).Insert()

# Shared constant values for these "Cebollín" records.
$mercadoId = 7
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region = "Ñuble"
$codreg = 16
$categoriaId = 100112037
$categoria = "Cebollín"
$variedad = "Sin especificar"
$clasificacion = "Hortaliza"

# New row 150
$r = 150
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 45106
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $categoriaId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 220
$ws.Cells.Item($r, 11).Value = 6000
$ws.Cells.Item($r, 12).Value = 7000
$ws.Cells.Item($r, 13).Value = 6545
$ws.Cells.Item($r, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item($r, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item($r, 16).Value = 182
$ws.Cells.Item($r, 17).Value = 36
$ws.Cells.Item($r, 18).Value = $clasificacion

# New row 151
$r = 151
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 45106
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $categoriaId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Segunda"
$ws.Cells.Item($r, 10).Value = 120
$ws.Cells.Item($r, 11).Value = 5000
$ws.Cells.Item($r, 12).Value = 5000
$ws.Cells.Item($r, 13).Value = 5000
$ws.Cells.Item($r, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item($r, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item($r, 16).Value = 139
$ws.Cells.Item($r, 17).Value = 36
$ws.Cells.Item($r, 18).Value = $clasificacion
